$d = $word.ActiveDocument

# Locate every occurrence of the literal run sequence "I öexj—Y" immediately
# followed by "–" and "I " (i.e. the 11-character span "I öexj—Y–I ") in
# document order. Only the 1st and 3rd occurrences (the two paragraphs that
# end shortly after with "Æj–¹sõ— ") are the ones touched by this revision;
# the 2nd/4th occurrences (inside the "öeZy…rçxI" paragraphs) are left
# untouched.

$needle = "I öexj—Y" + [char]8211 + "I "

$matches = New-Object System.Collections.ArrayList

$scan = $d.Content
$scan.Find.ClearFormatting()
$found = $scan.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $null = $matches.Add(@($scan.Start, $scan.End))
    $scan.Collapse(0)
    $scan.MoveEnd(1, 0)
    $scan.End = $d.Content.End
    $found = $scan.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

Write-Host "matches found:" $matches.Count

# Apply edits back-to-front so earlier offsets stay valid while we work.
for ($m = $matches.Count - 1; $m -ge 0; $m--) {
    $pair = $matches[$m]
    $matchStart = $pair[0]
    $matchEnd = $pair[1]

    # matchStart .. matchStart+8   -> run "I öexj—Y"      (8 chars)
    # matchStart+8 .. matchStart+9 -> run "–"              (1 char)
    # matchStart+9 .. matchEnd     -> run "I " (2 chars, matchEnd = matchStart+11)

    if ($m -eq 0) {
        # 1st occurrence -> split with a yellow-highlighted "Y", "–", "I"
        $rY = $d.Range($matchStart + 7, $matchStart + 8)
        $rY.Font.HighlightColorIndex = 7

        $rDash = $d.Range($matchStart + 8, $matchStart + 9)
        $rDash.Font.HighlightColorIndex = 7

        $rI = $d.Range($matchStart + 9, $matchStart + 10)
        $rI.Font.HighlightColorIndex = 7
    }
    elseif ($m -eq 2) {
        # 3rd occurrence -> drop the "–" and merge "Y" + "I" under a green highlight
        $rDash = $d.Range($matchStart + 8, $matchStart + 9)
        $rDash.Delete()

        # after deletion, "Y" sits at (matchStart+7, matchStart+8) and the
        # following "I" (previously at matchStart+9) has shifted left by 1
        $rYI = $d.Range($matchStart + 7, $matchStart + 9)
        $rYI.Font.HighlightColorIndex = 4
    }
    # other matches (2nd, 4th, ...) are left untouched
}

Write-Host "done"
